# Apply cell-level updates per the crypto price refresh diff.
# Values in column D that are numeric-looking strings (e.g. "217.42") are
# written with a leading apostrophe so Excel keeps them as literal Text
# (matching the source workbook, where every data cell is an inline string),
# instead of silently re-interpreting them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.208.64'
$ws.Range("E2").Value = '  -1.97%  '

# Row 3
$ws.Range("D3").Value = '1.669.80'
$ws.Range("E3").Value = '  -1.45%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").Value = '''217.42'
$ws.Range("E5").Value = '  -1.36%  '

# Row 6
$ws.Range("D6").Value = '''0.5111'
$ws.Range("E6").Value = '  -0.13%  '

# Row 7
$ws.Range("E7").Value = '  +0.07%  '

# Row 8
$ws.Range("D8").Value = '''0.2652'

# Row 9
$ws.Range("D9").Value = '''0.06369'
$ws.Range("E9").Value = '  +3.22%  '

# Row 10
$ws.Range("D10").Value = '''21.45'
$ws.Range("E10").Value = '  -1.99%  '

# Row 11
$ws.Range("D11").Value = '''0.07391'
$ws.Range("E11").Value = '  +0.71%  '

# Row 12
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = '''4.538'
$ws.Range("E12").Value = '  +1.70%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.677.77'
$ws.Range("E13").Value = '  -1.09%  '

# Row 14
$ws.Range("D14").Value = '''0.5824'
$ws.Range("E14").Value = '  +0.58%  '

# Row 15
$ws.Range("D15").Value = '''0.000008623'
$ws.Range("E15").Value = '  +5.98%  '

# Row 16
$ws.Range("D16").Value = '''64.43'
$ws.Range("E16").Value = '  -1.05%  '

# Row 17
$ws.Range("D17").Value = '26.260.17'
$ws.Range("E17").Value = '  -1.88%  '

# Row 18
$ws.Range("D18").Value = '''4.938'
$ws.Range("E18").Value = '  -1.00%  '

# Row 19
$ws.Range("D19").Value = '''1.006'
$ws.Range("E19").Value = '  +0.15%  '

# Row 20
$ws.Range("D20").Value = '''10.86'
$ws.Range("E20").Value = '  +2.16%  '

# Row 21
$ws.Range("D21").Value = '''188.94'
$ws.Range("E21").Value = '  +1.28%  '

# Row 22
$ws.Range("D22").Value = '''6.202'
$ws.Range("E22").Value = '  -0.79%  '

# Row 23
$ws.Range("E23").Value = '  -0.01%  '

# Row 24
$ws.Range("D24").Value = '''144.19'
$ws.Range("E24").Value = '  +1.21%  '

# Row 25
$ws.Range("D25").Value = '''7.638'
$ws.Range("E25").Value = '  +1.73%  '

# Row 26
$ws.Range("D26").Value = '''0.1177'
$ws.Range("E26").Value = '  +2.60%  '

# Row 27
$ws.Range("D27").Value = '''15.63'
$ws.Range("E27").Value = '  +2.91%  '

# Row 28
$ws.Range("D28").Value = '''0.05965'
$ws.Range("E28").Value = '  +1.54%  '

# Row 29
$ws.Range("D29").Value = '''1.285'
$ws.Range("E29").Value = '  -3.51%  '

# Row 30
$ws.Range("D30").Value = '''1.325'
$ws.Range("E30").Value = '  -1.59%  '

# Row 31
$ws.Range("D31").Value = '''3.520'
$ws.Range("E31").Value = '  +2.23%  '

# Row 32
$ws.Range("D32").Value = '''3.521'
$ws.Range("E32").Value = '  +2.99%  '

# Row 33
$ws.Range("D33").Value = '''1.644'
$ws.Range("E33").Value = '  +0.43%  '

# Row 34
$ws.Range("D34").Value = '''1.014'
$ws.Range("E34").Value = '  +2.90%  '

# Row 35
$ws.Range("D35").Value = '''0.6030'
$ws.Range("E35").Value = '  +0.89%  '

# Row 36
$ws.Range("D36").Value = '''2.376'
$ws.Range("E36").Value = '  -1.80%  '

# Row 37
$ws.Range("D37").Value = '''2.656'
$ws.Range("E37").Value = '  +0.08%  '

# Row 38
$ws.Range("B38").Value = 'FraxShare'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D38").Value = '''6.091'
$ws.Range("E38").Value = '  +4.13%  '

# Row 39
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '''0.01614'
$ws.Range("E39").Value = '  +1.33%  '

# Row 40
$ws.Range("D40").Value = '1.077.80'
$ws.Range("E40").Value = '  -0.81%  '

# Row 41
$ws.Range("D41").Value = '''0.8705'
$ws.Range("E41").Value = '  +1.04%  '

# Row 42
$ws.Range("E42").Value = '  +0.42%  '

# Row 43
$ws.Range("D43").Value = '''100.18'
$ws.Range("E43").Value = '  +2.74%  '

# Row 44
$ws.Range("D44").Value = '1.822.20'
$ws.Range("E44").Value = '  -1.12%  '

# Row 45
$ws.Range("E45").Value = '  +8.65%  '

# Row 46
$ws.Range("D46").Value = '''56.20'
$ws.Range("E46").Value = '  +0.56%  '

# Row 47
$ws.Range("E47").Value = '  +0.78%  '

# Row 48
$ws.Range("D48").Value = '''8.057'
$ws.Range("E48").Value = '  +1.79%  '

# Row 49
$ws.Range("D49").Value = '''0.05215'
$ws.Range("E49").Value = '  -0.46%  '

# Row 50
$ws.Range("D50").Value = '''0.4298'
$ws.Range("E50").Value = '  -0.36%  '

# Row 51
$ws.Range("D51").Value = '''5.883'
$ws.Range("E51").Value = '  +2.83%  '
